$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: reorder rows 3-6, 21-25, 32/34 (swap F:V content between rows) ---
# row 3 <= old row 6
$ws.Range("F3").Value = "Bayer Leverkusen"
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = "RB Leipzig"
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 2.59
$ws.Range("K3").Value = "01/07/2023 10:31"
$ws.Range("L3").Value = 2.48
$ws.Range("M3").Value = "19/08/2023 15:08"
$ws.Range("N3").Value = 3.38
$ws.Range("O3").Value = "01/07/2023 10:31"
$ws.Range("P3").Value = 3.63
$ws.Range("Q3").Value = "19/08/2023 15:29"
$ws.Range("R3").Value = 2.68
$ws.Range("S3").Value = "01/07/2023 10:31"
$ws.Range("T3").Value = 2.98
$ws.Range("U3").Value = "19/08/2023 15:29"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/germany/bundesliga/bayer-leverkusen-rb-leipzig/QNLpbj2b/"

# row 4 <= old row 3
$ws.Range("F4").Value = "Wolfsburg"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = "Heidenheim"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1.63
$ws.Range("K4").Value = "01/07/2023 10:29"
$ws.Range("L4").Value = 1.66
$ws.Range("M4").Value = "19/08/2023 15:14"
$ws.Range("N4").Value = 4.13
$ws.Range("O4").Value = "01/07/2023 10:29"
$ws.Range("P4").Value = 4.42
$ws.Range("Q4").Value = "19/08/2023 15:25"
$ws.Range("R4").Value = 4.88
$ws.Range("S4").Value = "01/07/2023 10:29"
$ws.Range("T4").Value = 5.16
$ws.Range("U4").Value = "19/08/2023 14:58"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/germany/bundesliga/wolfsburg-heidenheim/KMEgdUWA/"

# row 5 <= old row 4
$ws.Range("F5").Value = "Stuttgart"
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = "Bochum"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1.68
$ws.Range("K5").Value = "01/07/2023 10:27"
$ws.Range("L5").Value = 1.78
$ws.Range("M5").Value = "19/08/2023 15:25"
$ws.Range("N5").Value = 4.01
$ws.Range("O5").Value = "01/07/2023 10:27"
$ws.Range("P5").Value = 4.16
$ws.Range("Q5").Value = "19/08/2023 15:26"
$ws.Range("R5").Value = 4.61
$ws.Range("S5").Value = "01/07/2023 10:27"
$ws.Range("T5").Value = 4.52
$ws.Range("U5").Value = "19/08/2023 15:28"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/germany/bundesliga/vfb-stuttgart-bochum/AB76gSGT/"

# row 6 <= old row 5
$ws.Range("F6").Value = "Hoffenheim"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "Freiburg"
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 2.44
$ws.Range("K6").Value = "01/07/2023 10:29"
$ws.Range("L6").Value = 2.2
$ws.Range("M6").Value = "19/08/2023 15:29"
$ws.Range("N6").Value = 3.48
$ws.Range("O6").Value = "01/07/2023 10:29"
$ws.Range("P6").Value = 3.74
$ws.Range("Q6").Value = "19/08/2023 15:22"
$ws.Range("R6").Value = 2.8
$ws.Range("S6").Value = "01/07/2023 10:29"
$ws.Range("T6").Value = 3.32
$ws.Range("U6").Value = "19/08/2023 15:29"
$ws.Range("V6").Value = "https://www.betexplorer.com/football/germany/bundesliga/hoffenheim-freiburg/EqEcelnH/"

# row 21 <= old row 25
$ws.Range("F21").Value = "Werder Bremen"
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = "Mainz"
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 3.02
$ws.Range("K21").Value = "20/08/2023 09:02"
$ws.Range("L21").Value = 2.73
$ws.Range("M21").Value = "02/09/2023 15:24"
$ws.Range("N21").Value = 3.61
$ws.Range("O21").Value = "20/08/2023 09:02"
$ws.Range("P21").Value = 3.51
$ws.Range("Q21").Value = "02/09/2023 15:24"
$ws.Range("R21").Value = 2.23
$ws.Range("S21").Value = "20/08/2023 09:02"
$ws.Range("T21").Value = 2.7
$ws.Range("U21").Value = "02/09/2023 15:23"
$ws.Range("V21").Value = "https://www.betexplorer.com/football/germany/bundesliga/werder-bremen-mainz/0vVKTshq/"

# row 22 <= old row 23
$ws.Range("F22").Value = "Hoffenheim"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = "Wolfsburg"
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 2.56
$ws.Range("K22").Value = "20/08/2023 09:02"
$ws.Range("L22").Value = 2.58
$ws.Range("M22").Value = "02/09/2023 15:07"
$ws.Range("N22").Value = 3.59
$ws.Range("O22").Value = "20/08/2023 09:02"
$ws.Range("P22").Value = 3.97
$ws.Range("Q22").Value = "02/09/2023 15:28"
$ws.Range("R22").Value = 2.74
$ws.Range("S22").Value = "20/08/2023 09:02"
$ws.Range("T22").Value = 2.62
$ws.Range("U22").Value = "02/09/2023 15:12"
$ws.Range("V22").Value = "https://www.betexplorer.com/football/germany/bundesliga/hoffenheim-wolfsburg/fchhYNqS/"

# row 23 <= old row 24
$ws.Range("F23").Value = "Stuttgart"
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = "Freiburg"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1.82
$ws.Range("K23").Value = "20/08/2023 09:02"
$ws.Range("L23").Value = 2.35
$ws.Range("M23").Value = "02/09/2023 15:29"
$ws.Range("N23").Value = 3.92
$ws.Range("O23").Value = "20/08/2023 09:02"
$ws.Range("P23").Value = 3.64
$ws.Range("Q23").Value = "02/09/2023 15:29"
$ws.Range("R23").Value = 3.94
$ws.Range("S23").Value = "20/08/2023 09:02"
$ws.Range("T23").Value = 3.11
$ws.Range("U23").Value = "02/09/2023 15:29"
$ws.Range("V23").Value = "https://www.betexplorer.com/football/germany/bundesliga/vfb-stuttgart-freiburg/84YSRLMe/"

# row 24 <= old row 21
$ws.Range("F24").Value = "Augsburg"
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = "Bochum"
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = 1.92
$ws.Range("K24").Value = "20/08/2023 09:02"
$ws.Range("L24").Value = 2.25
$ws.Range("M24").Value = "02/09/2023 15:16"
$ws.Range("N24").Value = 3.74
$ws.Range("O24").Value = "20/08/2023 09:02"
$ws.Range("P24").Value = 3.68
$ws.Range("Q24").Value = "02/09/2023 15:16"
$ws.Range("R24").Value = 3.7
$ws.Range("S24").Value = "20/08/2023 09:02"
$ws.Range("T24").Value = 3.26
$ws.Range("U24").Value = "02/09/2023 15:16"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/germany/bundesliga/augsburg-bochum/vVTOS17k/"

# row 25 <= old row 22
$ws.Range("F25").Value = "Bayer Leverkusen"
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = "Darmstadt"
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 1.4
$ws.Range("K25").Value = "21/08/2023 06:14"
$ws.Range("L25").Value = 1.22
$ws.Range("M25").Value = "02/09/2023 15:24"
$ws.Range("N25").Value = 4.92
$ws.Range("O25").Value = "21/08/2023 06:14"
$ws.Range("P25").Value = 7.25
$ws.Range("Q25").Value = "02/09/2023 15:24"
$ws.Range("R25").Value = 6.76
$ws.Range("S25").Value = "21/08/2023 06:14"
$ws.Range("T25").Value = 13
$ws.Range("U25").Value = "02/09/2023 15:24"
$ws.Range("V25").Value = "https://www.betexplorer.com/football/germany/bundesliga/bayer-leverkusen-darmstadt/8tjtyQa9/"

# row 32 <= old row 34
$ws.Range("F32").Value = "Mainz"
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = "Stuttgart"
$ws.Range("I32").Value = 3
$ws.Range("J32").Value = 2.62
$ws.Range("K32").Value = "28/08/2023 16:01"
$ws.Range("L32").Value = 2.71
$ws.Range("M32").Value = "16/09/2023 15:25"
$ws.Range("N32").Value = 3.4
$ws.Range("O32").Value = "28/08/2023 16:01"
$ws.Range("P32").Value = 3.54
$ws.Range("Q32").Value = "16/09/2023 15:27"
$ws.Range("R32").Value = 2.79
$ws.Range("S32").Value = "28/08/2023 16:01"
$ws.Range("T32").Value = 2.71
$ws.Range("U32").Value = "16/09/2023 15:27"
$ws.Range("V32").Value = "https://www.betexplorer.com/football/germany/bundesliga/mainz-vfb-stuttgart/WMyiWzEs/"

# row 34 <= old row 32
$ws.Range("F34").Value = "FC Koln"
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = "Hoffenheim"
$ws.Range("I34").Value = 3
$ws.Range("J34").Value = 2.03
$ws.Range("K34").Value = "28/08/2023 16:01"
$ws.Range("L34").Value = 2.22
$ws.Range("M34").Value = "16/09/2023 15:17"
$ws.Range("N34").Value = 3.7
$ws.Range("O34").Value = "28/08/2023 16:01"
$ws.Range("P34").Value = 3.86
$ws.Range("Q34").Value = "16/09/2023 15:27"
$ws.Range("R34").Value = 3.68
$ws.Range("S34").Value = "28/08/2023 16:01"
$ws.Range("T34").Value = 3.19
$ws.Range("U34").Value = "16/09/2023 15:20"
$ws.Range("V34").Value = "https://www.betexplorer.com/football/germany/bundesliga/1-fc-koln-hoffenheim/lbWeVfTm/"

# --- Part 2: append new rows 39-46 ---
# new row 39
$ws.Range("A2:V2").Copy()
$ws.Range("A39:V39").PasteSpecial(-4122)
$ws.Range("B39").Value = "germany"
$ws.Range("C39").Value = "bundesliga"
$ws.Range("D39").Value = "2023-2024"
$ws.Range("A39").Value = 38
$ws.Range("E39").Value = 45192.64583333334
$ws.Range("F39").Value = "Dortmund"
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = "Wolfsburg"
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1.47
$ws.Range("K39").Value = "05/09/2023 12:01"
$ws.Range("L39").Value = 1.74
$ws.Range("M39").Value = "23/09/2023 15:20"
$ws.Range("N39").Value = 5.02
$ws.Range("O39").Value = "05/09/2023 12:01"
$ws.Range("P39").Value = 4.42
$ws.Range("Q39").Value = "23/09/2023 15:28"
$ws.Range("R39").Value = 5.4
$ws.Range("S39").Value = "05/09/2023 12:01"
$ws.Range("T39").Value = 4.48
$ws.Range("U39").Value = "23/09/2023 15:28"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/germany/bundesliga/dortmund-wolfsburg/dUMGQWsJ/"

# new row 40
$ws.Range("A2:V2").Copy()
$ws.Range("A40:V40").PasteSpecial(-4122)
$ws.Range("B40").Value = "germany"
$ws.Range("C40").Value = "bundesliga"
$ws.Range("D40").Value = "2023-2024"
$ws.Range("A40").Value = 39
$ws.Range("E40").Value = 45192.64583333334
$ws.Range("F40").Value = "B. Monchengladbach"
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = "RB Leipzig"
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = 4.05
$ws.Range("K40").Value = "05/09/2023 12:01"
$ws.Range("L40").Value = 4.91
$ws.Range("M40").Value = "23/09/2023 15:28"
$ws.Range("N40").Value = 4.16
$ws.Range("O40").Value = "05/09/2023 12:01"
$ws.Range("P40").Value = 4.55
$ws.Range("Q40").Value = "23/09/2023 15:28"
$ws.Range("R40").Value = 1.83
$ws.Range("S40").Value = "05/09/2023 12:01"
$ws.Range("T40").Value = 1.66
$ws.Range("U40").Value = "23/09/2023 15:28"
$ws.Range("V40").Value = "https://www.betexplorer.com/football/germany/bundesliga/b-monchengladbach-rb-leipzig/8M5YrEcm/"

# new row 41
$ws.Range("A2:V2").Copy()
$ws.Range("A41:V41").PasteSpecial(-4122)
$ws.Range("B41").Value = "germany"
$ws.Range("C41").Value = "bundesliga"
$ws.Range("D41").Value = "2023-2024"
$ws.Range("A41").Value = 40
$ws.Range("E41").Value = 45192.64583333334
$ws.Range("F41").Value = "Bayern Munich"
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = "Bochum"
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 1.07
$ws.Range("K41").Value = "05/09/2023 12:01"
$ws.Range("L41").Value = 1.12
$ws.Range("M41").Value = "23/09/2023 15:00"
$ws.Range("N41").Value = 20.79
$ws.Range("O41").Value = "05/09/2023 12:01"
$ws.Range("P41").Value = 10.5
$ws.Range("Q41").Value = "23/09/2023 14:59"
$ws.Range("R41").Value = 29.26
$ws.Range("S41").Value = "05/09/2023 12:01"
$ws.Range("T41").Value = 21
$ws.Range("U41").Value = "23/09/2023 15:00"
$ws.Range("V41").Value = "https://www.betexplorer.com/football/germany/bundesliga/bayern-munich-bochum/IwOCRCSC/"

# new row 42
$ws.Range("A2:V2").Copy()
$ws.Range("A42:V42").PasteSpecial(-4122)
$ws.Range("B42").Value = "germany"
$ws.Range("C42").Value = "bundesliga"
$ws.Range("D42").Value = "2023-2024"
$ws.Range("A42").Value = 41
$ws.Range("E42").Value = 45192.64583333334
$ws.Range("F42").Value = "Union Berlin"
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = "Hoffenheim"
$ws.Range("I42").Value = 2
$ws.Range("J42").Value = 1.85
$ws.Range("K42").Value = "05/09/2023 12:01"
$ws.Range("L42").Value = 2.13
$ws.Range("M42").Value = "23/09/2023 15:05"
$ws.Range("N42").Value = 3.93
$ws.Range("O42").Value = "05/09/2023 12:01"
$ws.Range("P42").Value = 3.72
$ws.Range("Q42").Value = "23/09/2023 15:20"
$ws.Range("R42").Value = 4.38
$ws.Range("S42").Value = "05/09/2023 12:01"
$ws.Range("T42").Value = 3.51
$ws.Range("U42").Value = "23/09/2023 15:26"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/germany/bundesliga/union-berlin-hoffenheim/AoQckGKJ/"

# new row 43
$ws.Range("A2:V2").Copy()
$ws.Range("A43:V43").PasteSpecial(-4122)
$ws.Range("B43").Value = "germany"
$ws.Range("C43").Value = "bundesliga"
$ws.Range("D43").Value = "2023-2024"
$ws.Range("A43").Value = 42
$ws.Range("E43").Value = 45192.64583333334
$ws.Range("F43").Value = "Augsburg"
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = "Mainz"
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = 2.6
$ws.Range("K43").Value = "05/09/2023 12:01"
$ws.Range("L43").Value = 2.56
$ws.Range("M43").Value = "23/09/2023 15:25"
$ws.Range("N43").Value = 3.57
$ws.Range("O43").Value = "05/09/2023 12:01"
$ws.Range("P43").Value = 3.52
$ws.Range("Q43").Value = "23/09/2023 15:05"
$ws.Range("R43").Value = 2.81
$ws.Range("S43").Value = "05/09/2023 12:01"
$ws.Range("T43").Value = 2.88
$ws.Range("U43").Value = "23/09/2023 15:25"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/germany/bundesliga/augsburg-mainz/tz4tshSa/"

# new row 44
$ws.Range("A2:V2").Copy()
$ws.Range("A44:V44").PasteSpecial(-4122)
$ws.Range("B44").Value = "germany"
$ws.Range("C44").Value = "bundesliga"
$ws.Range("D44").Value = "2023-2024"
$ws.Range("A44").Value = 43
$ws.Range("E44").Value = 45192.77083333334
$ws.Range("F44").Value = "Werder Bremen"
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = "FC Koln"
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 2.69
$ws.Range("K44").Value = "05/09/2023 12:01"
$ws.Range("L44").Value = 2.38
$ws.Range("M44").Value = "23/09/2023 18:28"
$ws.Range("N44").Value = 3.58
$ws.Range("O44").Value = "05/09/2023 12:01"
$ws.Range("P44").Value = 3.69
$ws.Range("Q44").Value = "23/09/2023 18:27"
$ws.Range("R44").Value = 2.69
$ws.Range("S44").Value = "05/09/2023 12:01"
$ws.Range("T44").Value = 3.02
$ws.Range("U44").Value = "23/09/2023 18:28"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/germany/bundesliga/werder-bremen-1-fc-koln/2q5xrYCg/"

# new row 45
$ws.Range("A2:V2").Copy()
$ws.Range("A45:V45").PasteSpecial(-4122)
$ws.Range("B45").Value = "germany"
$ws.Range("C45").Value = "bundesliga"
$ws.Range("D45").Value = "2023-2024"
$ws.Range("A45").Value = 44
$ws.Range("E45").Value = 45193.64583333334
$ws.Range("F45").Value = "Bayer Leverkusen"
$ws.Range("G45").Value = 4
$ws.Range("H45").Value = "Heidenheim"
$ws.Range("I45").Value = 1
$ws.Range("J45").Value = 1.29
$ws.Range("K45").Value = "11/09/2023 13:14"
$ws.Range("L45").Value = 1.2
$ws.Range("M45").Value = "24/09/2023 15:28"
$ws.Range("N45").Value = 5.81
$ws.Range("O45").Value = "11/09/2023 13:14"
$ws.Range("P45").Value = 7.91
$ws.Range("Q45").Value = "24/09/2023 15:28"
$ws.Range("R45").Value = 8.22
$ws.Range("S45").Value = "11/09/2023 13:14"
$ws.Range("T45").Value = 13.62
$ws.Range("U45").Value = "24/09/2023 15:29"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/germany/bundesliga/bayer-leverkusen-heidenheim/lxP1lzzQ/"

# new row 46
$ws.Range("A2:V2").Copy()
$ws.Range("A46:V46").PasteSpecial(-4122)
$ws.Range("B46").Value = "germany"
$ws.Range("C46").Value = "bundesliga"
$ws.Range("D46").Value = "2023-2024"
$ws.Range("A46").Value = 45
$ws.Range("E46").Value = 45193.72916666666
$ws.Range("F46").Value = "Eintracht Frankfurt"
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = "Freiburg"
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1.98
$ws.Range("K46").Value = "05/09/2023 12:01"
$ws.Range("L46").Value = 2.22
$ws.Range("M46").Value = "24/09/2023 17:22"
$ws.Range("N46").Value = 3.85
$ws.Range("O46").Value = "05/09/2023 12:01"
$ws.Range("P46").Value = 3.52
$ws.Range("Q46").Value = "24/09/2023 17:29"
$ws.Range("R46").Value = 3.87
$ws.Range("S46").Value = "05/09/2023 12:01"
$ws.Range("T46").Value = 3.47
$ws.Range("U46").Value = "24/09/2023 17:29"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/germany/bundesliga/eintracht-frankfurt-freiburg/fJ1Uqfrs/"

Write-Host "edit complete"
